$wb = $excel.ActiveWorkbook

# --- Sheet1 "Trends Status": Insufficient Data row totals 431 -> 432 ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B8").Value = 432
$ws1.Range("C8").Value = 432

# --- Sheet3 "Priority Status": update species counts ---
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# --- Sheet4 "Species qualification": rename assessment label + update count ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 432

# --- Rename existing "High Priority break-up" sheet and give it fresh data ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")
$ws5.Name = "Interannual update - High Pri"

# Add a brand-new sheet right after it that keeps the old break-up data
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "Major update - High Priority "

# Populate the new "Major update" sheet with the previous break-up contents
$ws6.Range("A1").Value = "Break-up"
$ws6.Range("B1").Value = "High Species (no.)"
$ws6.Range("C1").Value = "High Species (perc.)"
$ws6.Range("D1").Value = "New High Species (no.)"
$ws6.Range("E1").Value = "New High Species (perc.)"
$ws6.Range("A1:E1").Font.Bold = $true
$ws6.Range("A1:E1").HorizontalAlignment = -4108

$ws6.Range("A2").Value = "IUCN"
$ws6.Range("B2").Value = 18
$ws6.Range("C2").Value = 100
$ws6.Range("D2").Value = 18
$ws6.Range("E2").Value = 100

# Replace the renamed sheet's rows with the interannual-update figures
$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 79
$ws5.Range("C2").Value = 76.7
$ws5.Range("D2").Value = 79
$ws5.Range("E2").Value = 89.8

$ws5.Range("A3").Value = "IUCN"
$ws5.Range("B3").Value = 24
$ws5.Range("C3").Value = 23.3
$ws5.Range("D3").Value = 9
$ws5.Range("E3").Value = 10.2

# Keep the originally-selected tab ("Trends Status") active, as before
$ws1.Activate()
